# Adds row 23 to the "Artfynd" sheet (new species observation record),
# matching the appended <row r="23"> block from the source diff.
# Dimension ref will be recalculated automatically by the engine on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

# --- Numeric cells -----------------------------------------------------
$ws.Cells.Item($row, 1).Value  = 131219383   # A  Id
$ws.Cells.Item($row, 2).Value  = 57881       # B  Taxonsorteringsordning
$ws.Cells.Item($row, 5).Value  = 100049      # E  TaxonId
$ws.Cells.Item($row, 17).Value = 516893      # Q  Ost
$ws.Cells.Item($row, 18).Value = 6705745     # R  Nord
$ws.Cells.Item($row, 19).Value = 25          # S  Noggrannhet

# --- Plain text cells ---------------------------------------------------
$ws.Cells.Item($row, 4).Value  = "NT"                                              # D  Rodlistade
$ws.Cells.Item($row, 6).Value  = "Spillkråka"                                      # F  Artnamn
$ws.Cells.Item($row, 7).Value  = "Dryocopus martius"                               # G  Vetenskapligt namn
$ws.Cells.Item($row, 8).Value  = "(Linnaeus, 1758)"                                # H  Auktor
$ws.Cells.Item($row, 13).Value = "gammalt bo"                                      # M  Aktivitet
$ws.Cells.Item($row, 16).Value = "Bjusan, Dlr"                                     # P  Lokalnamn
$ws.Cells.Item($row, 20).Value = "Dalarna"                                         # T  Lan
$ws.Cells.Item($row, 21).Value = "Borlänge"                                        # U  Kommun
$ws.Cells.Item($row, 22).Value = "Dalarna"                                         # V  Provins
$ws.Cells.Item($row, 23).Value = "Stora Tuna"                                      # W  Socken
$ws.Cells.Item($row, 29).Value = "Grov äldre tall. 10 meter upp, ovalt hål, cirka 7 cm stort."  # AC Publik kommentar
$ws.Cells.Item($row, 49).Value = "Anna-Lena Thommson"                              # AW Rapportor
$ws.Cells.Item($row, 50).Value = "Anna-Lena Thommson"                              # AX Observatorer

# Date-looking text must stay text (not be auto-converted to a date serial),
# so force it with a leading apostrophe (Excel's "enter as text" convention).
$ws.Cells.Item($row, 25).Value = "'2026-02-14"   # Y  Startdatum
$ws.Cells.Item($row, 27).Value = "'2026-02-14"   # AA Slutdatum

# --- Boolean cells -------------------------------------------------------
$ws.Cells.Item($row, 30).Value = $false   # AD Ej aterfunnen
$ws.Cells.Item($row, 31).Value = $false   # AE Osaker artbestamning
$ws.Cells.Item($row, 33).Value = $false   # AG Ospontan

# --- Present-but-empty text cells ---------------------------------------
# A bare "" clears/omits the cell in this engine, so use the apostrophe
# trick to force a genuine (empty) text cell to exist, matching the
# self-closed <c t="inlineStr"/> cells from the source row.
$ws.Cells.Item($row, 9).Value  = "'"   # I  Antal
$ws.Cells.Item($row, 11).Value = "'"   # K  Alder-Stadium
$ws.Cells.Item($row, 12).Value = "'"   # L  Kon
$ws.Cells.Item($row, 14).Value = "'"   # N  Metod
$ws.Cells.Item($row, 46).Value = "'"   # AT Bestamningsar
$ws.Cells.Item($row, 51).Value = "'"   # AY Projektnamn
